# "Assets and Tracking Script added"
# - Add a new "Company_Profile" worksheet at the end of the workbook with
#   Name / DOT / Email rows.
# - Flip Configuration!Send Mail from "yes" to "no".

$wb = $excel.ActiveWorkbook

# --- 1. Add the new Company_Profile sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Company_Profile"

$newSheet.Range("A1").Value = "Name"
$newSheet.Range("A2").Value = "DOT"
$newSheet.Range("A3").Value = "Email"

$newSheet.Range("B1").Value = "Zeus Transport Inc. 3"
$newSheet.Range("B2").Value = 123412346
$newSheet.Range("B3").Value = "demo@truckx.com"

# Text format (matches the "Text" style used by the other input sheets)
$newSheet.Range("A1:B3").NumberFormat = "@"

# Column widths approximating the source workbook's layout
$newSheet.Columns.Item(1).ColumnWidth = 14.79
$newSheet.Columns.Item(2).ColumnWidth = 20.5

# --- 2. Update Configuration sheet: "Send Mail" -> "no" ---
$config = $wb.Worksheets.Item("Configuration")
$config.Range("B3").Value = "no"
